$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 63, pushing the existing rows 63-80 down to 64-81.
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row with this week's price record.
$ws.Range("A63").Value = 9
$ws.Range("B63").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C63").Value = "Metropolitana"
$ws.Range("D63").Value = 44785
$ws.Range("E63").Value = 13
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100102
$ws.Range("H63").Value = "Cítricos"
$ws.Range("I63").Value = 100102006
$ws.Range("J63").Value = "Pomelo"
$ws.Range("K63").Value = "Start Ruby"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 400
$ws.Range("N63").Value = 12000
$ws.Range("O63").Value = 12000
$ws.Range("P63").Value = 12000
$ws.Range("Q63").Value = "$/caja 14 kilos"
$ws.Range("R63").Value = "Región Metropolitana"
$ws.Range("S63").Value = 857
$ws.Range("T63").Value = 14
